$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 84 by copying row 83 formatting first, then set its values below
$ws.Range("A83:V83").Copy($ws.Range("A84:V84"))

# Row 52 <= data from old row 53
$ws.Range("F52").Value = "Cibalia"
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = "Bijelo Brdo"
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 1.69
$ws.Range("K52").Value = "06/10/2023 02:12"
$ws.Range("L52").Value = 1.78
$ws.Range("M52").Value = "07/10/2023 14:42"
$ws.Range("N52").Value = 3.44
$ws.Range("O52").Value = "06/10/2023 02:12"
$ws.Range("P52").Value = 3.57
$ws.Range("Q52").Value = "07/10/2023 14:42"
$ws.Range("R52").Value = 4.27
$ws.Range("S52").Value = "06/10/2023 02:12"
$ws.Range("T52").Value = 4.38
$ws.Range("U52").Value = "07/10/2023 14:42"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-bijelo-brdo/IVBVDA5n/"

# Row 53 <= data from old row 52
$ws.Range("F53").Value = "Dugopolje"
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = "Zrinski Jurjevac"
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 2.05
$ws.Range("K53").Value = "06/10/2023 02:12"
$ws.Range("L53").Value = 2.15
$ws.Range("M53").Value = "07/10/2023 14:02"
$ws.Range("N53").Value = 3.21
$ws.Range("O53").Value = "06/10/2023 02:12"
$ws.Range("P53").Value = 3.39
$ws.Range("Q53").Value = "07/10/2023 14:02"
$ws.Range("R53").Value = 3.17
$ws.Range("S53").Value = "06/10/2023 02:12"
$ws.Range("T53").Value = 3.2
$ws.Range("U53").Value = "07/10/2023 14:02"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-zrinski-jurjevac/Ec5vClza/"

# Row 64 <= data from old row 65
$ws.Range("F64").Value = "Vukovar 1991"
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = "Solin"
$ws.Range("I64").Value = 2
$ws.Range("J64").Value = 1.5
$ws.Range("K64").Value = "20/10/2023 02:12"
$ws.Range("L64").Value = 1.6
$ws.Range("M64").Value = "21/10/2023 14:58"
$ws.Range("N64").Value = 4.02
$ws.Range("O64").Value = "20/10/2023 02:12"
$ws.Range("P64").Value = 4.19
$ws.Range("Q64").Value = "21/10/2023 14:58"
$ws.Range("R64").Value = 4.93
$ws.Range("S64").Value = "20/10/2023 02:12"
$ws.Range("T64").Value = 4.78
$ws.Range("U64").Value = "21/10/2023 14:56"
$ws.Range("V64").Value = "https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-solin/jHPc3Hmd/"

# Row 65 <= data from old row 64
$ws.Range("F65").Value = "Dugopolje"
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = "Cibalia"
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1.94
$ws.Range("K65").Value = "20/10/2023 02:12"
$ws.Range("L65").Value = 1.71
$ws.Range("M65").Value = "21/10/2023 10:12"
$ws.Range("N65").Value = 3.21
$ws.Range("O65").Value = "20/10/2023 02:12"
$ws.Range("P65").Value = 3.45
$ws.Range("Q65").Value = "21/10/2023 13:03"
$ws.Range("R65").Value = 3.45
$ws.Range("S65").Value = "20/10/2023 02:12"
$ws.Range("T65").Value = 5.04
$ws.Range("U65").Value = "21/10/2023 14:09"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-cibalia/QFsPNoZA/"

# Row 76 <= data from old row 78
$ws.Range("F76").Value = "Orijent"
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = "Jarun"
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 1.79
$ws.Range("K76").Value = "03/11/2023 02:12"
$ws.Range("L76").Value = 1.79
$ws.Range("M76").Value = "04/11/2023 13:52"
$ws.Range("N76").Value = 3.62
$ws.Range("O76").Value = "03/11/2023 02:12"
$ws.Range("P76").Value = 4.01
$ws.Range("Q76").Value = "04/11/2023 13:52"
$ws.Range("R76").Value = 3.55
$ws.Range("S76").Value = "03/11/2023 02:12"
$ws.Range("T76").Value = 3.79
$ws.Range("U76").Value = "04/11/2023 13:52"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-jarun/QZwtGfX1/"

# Row 77 <= data from old row 76
$ws.Range("F77").Value = "Cibalia"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Dubrava"
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1.98
$ws.Range("K77").Value = "03/11/2023 02:12"
$ws.Range("L77").Value = 2.56
$ws.Range("M77").Value = "04/11/2023 13:59"
$ws.Range("N77").Value = 3.32
$ws.Range("O77").Value = "03/11/2023 02:12"
$ws.Range("P77").Value = 3.28
$ws.Range("Q77").Value = "04/11/2023 13:59"
$ws.Range("R77").Value = 3.25
$ws.Range("S77").Value = "03/11/2023 02:12"
$ws.Range("T77").Value = 2.67
$ws.Range("U77").Value = "04/11/2023 13:59"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-dubrava-zagreb/nquxHzIe/"

# Row 78 <= data from old row 77
$ws.Range("F78").Value = "Dugopolje"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Bijelo Brdo"
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 1.62
$ws.Range("K78").Value = "03/11/2023 02:12"
$ws.Range("L78").Value = 1.92
$ws.Range("M78").Value = "04/11/2023 13:59"
$ws.Range("N78").Value = 3.57
$ws.Range("O78").Value = "03/11/2023 02:12"
$ws.Range("P78").Value = 3.03
$ws.Range("Q78").Value = "04/11/2023 13:59"
$ws.Range("R78").Value = 4.6
$ws.Range("S78").Value = "03/11/2023 02:12"
$ws.Range("T78").Value = 4.54
$ws.Range("U78").Value = "04/11/2023 13:59"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-bijelo-brdo/jHoTIdmq/"

# Row 82 <= data from old row 83
$ws.Range("F82").Value = "Vukovar 1991"
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = "Orijent"
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1.47
$ws.Range("K82").Value = "10/11/2023 02:12"
$ws.Range("L82").Value = 1.53
$ws.Range("M82").Value = "11/11/2023 13:54"
$ws.Range("N82").Value = 4.22
$ws.Range("O82").Value = "10/11/2023 02:12"
$ws.Range("P82").Value = 4.56
$ws.Range("Q82").Value = "11/11/2023 13:54"
$ws.Range("R82").Value = 5.03
$ws.Range("S82").Value = "10/11/2023 02:12"
$ws.Range("T82").Value = 5
$ws.Range("U82").Value = "11/11/2023 13:54"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-orijent/8WYwyd9R/"

# Row 83 <= data from old row 82
$ws.Range("F83").Value = "Bijelo Brdo"
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = "Sibenik"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 4.54
$ws.Range("K83").Value = "10/11/2023 02:12"
$ws.Range("L83").Value = 5.86
$ws.Range("M83").Value = "11/11/2023 13:52"
$ws.Range("N83").Value = 3.58
$ws.Range("O83").Value = "10/11/2023 02:12"
$ws.Range("P83").Value = 4.01
$ws.Range("Q83").Value = "11/11/2023 13:52"
$ws.Range("R83").Value = 1.63
$ws.Range("S83").Value = "10/11/2023 02:12"
$ws.Range("T83").Value = 1.54
$ws.Range("U83").Value = "11/11/2023 13:51"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-sibenik/z5LWxIvF/"

# Set new row 84 content (betexplorer match: Croatia Zmijavci vs Dugopolje)
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "croatia"
$ws.Range("C84").Value = "prva-nl"
$ws.Range("D84").Value = "2023-2024"
$ws.Range("E84").Value = 45242.58333333334
$ws.Range("F84").Value = "Croatia Zmijavci"
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = "Dugopolje"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 2.31
$ws.Range("K84").Value = "11/11/2023 02:12"
$ws.Range("L84").Value = 2.42
$ws.Range("M84").Value = "12/11/2023 13:43"
$ws.Range("N84").Value = 3.11
$ws.Range("O84").Value = "11/11/2023 02:12"
$ws.Range("P84").Value = 2.96
$ws.Range("Q84").Value = "12/11/2023 13:43"
$ws.Range("R84").Value = 2.85
$ws.Range("S84").Value = "11/11/2023 02:12"
$ws.Range("T84").Value = 3.12
$ws.Range("U84").Value = "12/11/2023 13:43"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/croatia/prva-nl/croatia-zmijavci-dugopolje/67EEUgne/"
